# This workbook tracks daily Lemon ("Limón") prices for "Feria Lagunitas de
# Puerto Montt". A new pair of price observations (dated 2021-11-05, serial
# 44505) is being inserted at the top of the price history block (rows
# 267-268), pushing all the later (older) observations down by two rows.
# The sheet dimension grows from A1:T306 to A1:T308.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first price row of the block
# (row 267), shifting every existing row at/after 267 down by two.
$ws.Range("A267:A268").EntireRow.Insert()

# --- New row 267: "1a amarillo" ---
$ws.Cells.Item(267, 1).Value = 4
$ws.Cells.Item(267, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(267, 3).Value = "Los Lagos"
$ws.Cells.Item(267, 4).Value = 44505
$ws.Cells.Item(267, 5).Value = 10
$ws.Cells.Item(267, 6).Value = "Fruta"
$ws.Cells.Item(267, 7).Value = 100102
$ws.Cells.Item(267, 8).Value = "Cítricos"
$ws.Cells.Item(267, 9).Value = 100102003
$ws.Cells.Item(267, 10).Value = "Limón"
$ws.Cells.Item(267, 11).Value = "Sin especificar"
$ws.Cells.Item(267, 12).Value = "1a amarillo"
$ws.Cells.Item(267, 13).Value = 1300
$ws.Cells.Item(267, 14).Value = 11000
$ws.Cells.Item(267, 15).Value = 12000
$ws.Cells.Item(267, 16).Value = 11500
$ws.Cells.Item(267, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(267, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(267, 19).Value = 639
$ws.Cells.Item(267, 20).Value = 18

# --- New row 268: "2a amarillo" ---
$ws.Cells.Item(268, 1).Value = 4
$ws.Cells.Item(268, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(268, 3).Value = "Los Lagos"
$ws.Cells.Item(268, 4).Value = 44505
$ws.Cells.Item(268, 5).Value = 10
$ws.Cells.Item(268, 6).Value = "Fruta"
$ws.Cells.Item(268, 7).Value = 100102
$ws.Cells.Item(268, 8).Value = "Cítricos"
$ws.Cells.Item(268, 9).Value = 100102003
$ws.Cells.Item(268, 10).Value = "Limón"
$ws.Cells.Item(268, 11).Value = "Sin especificar"
$ws.Cells.Item(268, 12).Value = "2a amarillo"
$ws.Cells.Item(268, 13).Value = 300
$ws.Cells.Item(268, 14).Value = 9000
$ws.Cells.Item(268, 15).Value = 9000
$ws.Cells.Item(268, 16).Value = 9000
$ws.Cells.Item(268, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(268, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(268, 19).Value = 500
$ws.Cells.Item(268, 20).Value = 18
